$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 1 (the "2021" header row), shifting all subsequent rows up by one.
$ws.Rows("1").Delete()

# Excel's row-delete normally re-points any defined name / print area that
# referenced the deleted row automatically; make sure all three Print_Area
# style names end up pointing one row higher, matching the new layout.
$sheetName = $ws.Name
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Print_Area") {
        $n.RefersTo = "='" + $sheetName + "'!`$A`$1:`$F`$20"
    } elseif ($n.Name -like "*!Print_Area") {
        $n.RefersTo = "='" + $sheetName + "'!`$A`$1:`$U`$19"
    } elseif ($n.Name -like "*!Print_Area_MI") {
        $n.RefersTo = "='" + $sheetName + "'!`$A`$1:`$F`$20"
    }
}

# Match the final selection recorded in the saved file.
$ws.Range("B21").Select()
